# Actualización automática del tracker
# Adds new result rows (112-120) to the tracker_resultados sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records to append at the bottom of the table.
$newRows = @(
    @{ A = 14831274; B = "2025-10-07"; C = "Dominique Rolland"; D = "Maxwell McKennon";      E = "Gana Dominique Rolland";         F = 1.83 },
    @{ A = 14841551; B = "2025-10-07"; C = "Olle Wallin";       D = "Edward Winter";          E = "Gana Edward Winter";             F = 2.1  },
    @{ A = 14831276; B = "2025-10-07"; C = "Benjamin Hassan";   D = "Tibo Colson";            E = "Gana Tibo Colson";               F = 3.25 },
    @{ A = 14841553; B = "2025-10-07"; C = "Marek Gengel";      D = "Evan Bynoe";             E = "Gana Evan Bynoe";                F = 2.75 },
    @{ A = 14841554; B = "2025-10-07"; C = "Aryan Shah";        D = "Jay Dylan Hara Friend";  E = "Gana Jay Dylan Hara Friend";     F = 1.62 },
    @{ A = 14831272; B = "2025-10-07"; C = "Bor Artnak";        D = "Abdullah Shelbayh";      E = "Gana Bor Artnak";                F = 2.75 },
    @{ A = 14841552; B = "2025-10-07"; C = "Daniel Milavsky";   D = "Darian King";            E = "Gana Darian King";               F = 2.63 },
    @{ A = 14841555; B = "2025-10-07"; C = "Michael Mmoh";      D = "Quinn Vandecasteele";    E = "Gana Quinn Vandecasteele";       F = 3    },
    @{ A = 14838649; B = "2025-10-07"; C = "Hynek Barton";      D = "Max Basing";             E = "Gana Hynek Barton";              F = 2    }
)

$startRow = 112

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A

    # Force column B to stay as literal text (e.g. "2025-10-07") instead of
    # being auto-converted into a date serial number.
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $row.B

    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    # resultado (G) and profit (H) are left blank - the match hasn't been
    # settled yet, matching the other pending rows already in the sheet.
}

Write-Host "Added $($newRows.Count) rows starting at row $startRow"
